# The author's export tool re-sorted the "Artfynd" observation rows: the
# same 9 observation records (rows 2..10) got shuffled onto different
# physical rows while row 1 (the header) stays put. This performs that
# single 9-cycle permutation of rows 2..10, preserving each record's full
# set of column values/types exactly (numbers stay numbers, booleans stay
# booleans, and text that merely looks like a number/date/time - e.g. the
# "Antal" count column or the Start/Sluttid date & time columns, all stored
# as text in the source file - is kept as text rather than being
# reinterpreted by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 10
$lastCol  = 51   # column AY

# New row -> old row it now contains (a single 9-cycle: 2<-9<-10<-5<-8<-4<-7<-6<-3<-2).
$mapping = @{
    2  = 9
    3  = 2
    4  = 7
    5  = 8
    6  = 3
    7  = 6
    8  = 4
    9  = 10
    10 = 5
}

# Columns whose text content could be misread as a number/date/time by Excel's
# normal type-inference on write (e.g. "11", "2023-09-03", "10:49"), even
# though the source file stores them as plain text. Force Text format on
# these before writing so the value round-trips as text, not a number/date.
$textForceCols = @(9, 25, 26, 27, 28)  # I, Y, Z, AA, AB

# 1) Snapshot every source row's values (and, implicitly, types) before any
#    writes happen - the permutation is a single 9-cycle so every row must be
#    captured before it gets overwritten.
$snapshots = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rng = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $lastCol))
    $snapshots[$r] = $rng.Value2
}

# 2) Pre-format the risky columns across the whole block as Text so the
#    upcoming writes don't get reinterpreted as numbers/dates/times.
foreach ($col in $textForceCols) {
    $colRange = $ws.Range($ws.Cells.Item($firstRow, $col), $ws.Cells.Item($lastRow, $col))
    $colRange.NumberFormat = "@"
}

# 3) Write each destination row from its recorded source snapshot.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $destRange = $ws.Range($ws.Cells.Item($destRow, 1), $ws.Cells.Item($destRow, $lastCol))
    $destRange.Value2 = $snapshots[$srcRow]
}

# 4) Drop the temporary Text formatting again so the cells end up with the
#    same (default) style they had before, now that the values are committed
#    as text.
foreach ($col in $textForceCols) {
    $colRange = $ws.Range($ws.Cells.Item($firstRow, $col), $ws.Cells.Item($lastRow, $col))
    $colRange.ClearFormats()
}
